# This deck ships two theme parts:
#   ppt/theme/theme1.xml -> currently the stock "Office Theme" colour scheme
#                            (used only by the Notes Master)
#   ppt/theme/theme2.xml -> currently the "Integral" colour scheme
#                            (used by the Slide Master / the whole deck)
#
# The authored change swaps the content of the two theme parts: the
# palette that used to live in theme2.xml ("Integral") moves to
# theme1.xml, and the stock "Office Theme" palette that used to live in
# theme1.xml moves to theme2.xml.
#
# The PowerPoint object model only exposes one live ColorScheme for the
# whole presentation (reached through the Slide Master / any slide), and
# it always maps onto ppt/theme/theme2.xml - the theme actually applied
# to the deck. We therefore re-point that scheme at the "Office Theme"
# palette, which is the reachable half of the swap (and the half that
# actually affects how the deck renders).

function HexToBGR($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p  = $ppt.ActivePresentation
$m  = $p.SlideMaster
$cs = $m.ColorScheme

# Index -> theme slot:
#  1 dk1   2 lt1   3 dk2     4 lt2
#  5 accent1  6 accent2  7 accent3  8 accent4
#  9 accent5 10 accent6 11 hlink   12 folHlink
$cs.Colors(1).RGB  = HexToBGR "000000"   # dk1
$cs.Colors(2).RGB  = HexToBGR "FFFFFF"   # lt1
$cs.Colors(3).RGB  = HexToBGR "44546A"   # dk2
$cs.Colors(4).RGB  = HexToBGR "E7E6E6"   # lt2
$cs.Colors(5).RGB  = HexToBGR "5B9BD5"   # accent1
$cs.Colors(6).RGB  = HexToBGR "ED7D31"   # accent2
$cs.Colors(7).RGB  = HexToBGR "A5A5A5"   # accent3
$cs.Colors(8).RGB  = HexToBGR "FFC000"   # accent4
$cs.Colors(9).RGB  = HexToBGR "4472C4"   # accent5
$cs.Colors(10).RGB = HexToBGR "70AD47"   # accent6
$cs.Colors(11).RGB = HexToBGR "0563C1"   # hlink
$cs.Colors(12).RGB = HexToBGR "954F72"   # folHlink
